$wb = $excel.ActiveWorkbook

# --- Transactions sheet: remove the two trade rows (rows 2 and 3) ---
$ws = $wb.Worksheets.Item("Transactions")
$ws.Rows("2:3").Delete()

# Reset the AutoFilter so its range shrinks to the remaining header-only data (A1:M1)
if ($ws.AutoFilterMode) {
    $ws.AutoFilterMode = $false
}
$ws.Range("A1:M1").AutoFilter(1)

# Update the hidden _FilterDatabase defined name to match the new range
foreach ($n in $wb.Names) {
    if ($n.Name -like "*_FilterDatabase*") {
        $n.RefersTo = "='Transactions'!`$A`$1:`$M`$1"
    }
}

# --- Summary sheet: zero-out metrics that depended on the removed trades ---
$ws2 = $wb.Worksheets.Item("Summary")
$ws2.Range("B2").Value = 0   # Total Trades (entries+exits)
$ws2.Range("B3").Value = 0   # Total PnL
$ws2.Range("B4").Value = 0   # Profit target hits
